$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "visual stimulus (IAPS) appeared for 2000 msec." -> "... 5000 msec."
#    (only this one occurrence changes; the other two "2000 msec" mentions
#    in the same paragraph stay untouched)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("(IAPS) appeared for 2000 msec", $true, $false, $false, `
    $false, $false, $true, 1, $false, "(IAPS) appeared for 5000 msec", 2)

# ---------------------------------------------------------------------------
# 2. Relocate the long block of text that currently sits right after the
#    "_GoBack" bookmark (Participants were seated ... through ... Pupil size
#    and gaze position were measured during the entirety of the trial.) so
#    that it instead sits right BEFORE the bookmark - i.e. swap the block
#    with the (empty) bookmark that used to precede it.
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bmStart = $bm.Start

$tail = $d.Range($bmStart, $d.Content.End)
$tail.Find.Execute("Pupil size and gaze position were measured during the entirety of the trial. ", `
    $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$blockEnd = $tail.End

$block = $d.Range($bmStart, $blockEnd)
$block.Cut()

$bm2 = $d.Bookmarks("_GoBack")
$dest = $d.Range($bm2.Start, $bm2.Start)
$dest.Paste()

# ---------------------------------------------------------------------------
# 3. The heading "References" no longer falls at the start of a printed
#    page once the block above moves earlier, so the stale
#    lastRenderedPageBreak marker in front of it has to go. The Word object
#    model has no direct handle for that marker, so the heading run is
#    deleted and retyped (preserving its Bold formatting), which drops it.
# ---------------------------------------------------------------------------
$heading = $d.Content
$heading.Find.Execute("References", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$headingStart = $heading.Start
$heading.Delete()

$newHeading = $d.Range($headingStart, $headingStart)
$newHeading.InsertBefore("References")
$newHeadingRange = $d.Range($headingStart, $headingStart + 10)
$newHeadingRange.Bold = 1
